$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91: "Clase 06" section header ---
$ws.Range("A91").Value = "Clase 06"
$ws.Range("A8:B8").Copy() | Out-Null
$ws.Range("A91:B91").PasteSpecial(-4122) | Out-Null

# --- Rows 92-103: entered sequentially, in order ---
$ws.Range("A92").Value = 0.00312500000000000017
$ws.Range("B92").Value = "Protocolo HTTP"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A92:B92").PasteSpecial(-4122) | Out-Null

$ws.Range("A93").Value = 0.00972222222222222238
$ws.Range("B93").Value = "cliente => request  // servidor => responses   ( Frontend // Backend )"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A93:B93").PasteSpecial(-4122) | Out-Null

$ws.Range("A94").Value = 0.01388888888888888812
$ws.Range("B94").Value = "Instalar nodemon de manera global (ver fila 84)"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A94:B94").PasteSpecial(-4122) | Out-Null

$ws.Range("A95").Value = 0.01631944444444444545
$ws.Range("B95").Value = "Propiedad `"type`": `"module`" para activar import y export "
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A95:B95").PasteSpecial(-4122) | Out-Null

$ws.Range("A96").Value = 0.01909722222222222376
$ws.Range("B96").Value = "script `"dev`": `"nodemon server-http .js`""
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A96:B96").PasteSpecial(-4122) | Out-Null

$ws.Range("A97").Value = 0.02152777777777777762
$ws.Range("B97").Value = "request (Propiedades, solicitado por el cliente) y response en el createServer"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A97:B97").PasteSpecial(-4122) | Out-Null

$ws.Range("A98").Value = 0.02326388888888888951
$ws.Range("B98").Value = "res.end enviar informacion al cliente"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A98:B98").PasteSpecial(-4122) | Out-Null

$ws.Range("A99").Value = 0.02430555555555555594
$ws.Range("B99").Value = "Thunder Client (Cliente https)"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A99:B99").PasteSpecial(-4122) | Out-Null

$ws.Range("A100").Value = 0.02534722222222222238
$ws.Range("B100").Value = "Postman, el que se usara en clase"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A100:B100").PasteSpecial(-4122) | Out-Null

$ws.Range("A101").Value = 0.02812500000000000069
$ws.Range("B101").Value = "Endpoints"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A101:B101").PasteSpecial(-4122) | Out-Null

$ws.Range("A102").Value = 0.03055555555555555455
$ws.Range("B102").Value = "Creacion de un Endpoint"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A102:B102").PasteSpecial(-4122) | Out-Null

$ws.Range("A103").Value = 0.03680555555555555663
$ws.Range("B103").Value = "Objeto Request  ( req )"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A103:B103").PasteSpecial(-4122) | Out-Null

# --- Rows 104/105: "https://expressjs.com" typed first, then a row was
# inserted above it for "Express - npm i express" (reproduces the shared-
# string index order: B105 -> idx 103, B104 -> idx 104) ---
$ws.Range("A104").Value = 0.04444444444444444614
$ws.Range("B104").Value = "https://expressjs.com"
$ws.Range("A21:B21").Copy() | Out-Null
$ws.Range("A104:B104").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(104).Insert() | Out-Null
$ws.Range("A104").Value = 0.04374999999999999722
$ws.Range("B104").Value = "Express  -  npm i express"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A104:B104").PasteSpecial(-4122) | Out-Null

# --- Rows 106-110: rows for "Express aplica...", "Desde el servidor...",
# "Propiedad req.query", "Propiedad req.params" were entered first forming
# consecutive rows 106-109; then a row was inserted at 106 for "Metodo get"
# (reproduces shared-string index order: those 4 get idx 105-108, "Metodo
# get" gets idx 109) ---
$ws.Range("A106").Value = 0.05034722222222222376
$ws.Range("B106").Value = "Express aplica el tipo de archivo sin usar el stringify, devuelve los datos en formato JSON"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A106:B106").PasteSpecial(-4122) | Out-Null

$ws.Range("A107").Value = 0.05138888888888888673
$ws.Range("B107").Value = "Desde el servidor somos los responsables de crear los estatus, errores del servidor"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A107:B107").PasteSpecial(-4122) | Out-Null

$ws.Range("A108").Value = 0.06006944444444444614
$ws.Range("B108").Value = "Propiedad req.query"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A108:B108").PasteSpecial(-4122) | Out-Null

$ws.Range("A109").Value = 0.07187499999999999445
$ws.Range("B109").Value = "Propiedad req.params"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A109:B109").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(106).Insert() | Out-Null
$ws.Range("A106").Value = 0.04722222222222222099
$ws.Range("B106").Value = "Metodo get"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A106:B106").PasteSpecial(-4122) | Out-Null

# --- Rows 111-112: entered sequentially, in order ---
$ws.Range("A111").Value = 0.07708333333333333703
$ws.Range("B111").Value = "Metodo post"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A111:B111").PasteSpecial(-4122) | Out-Null

$ws.Range("A112").Value = 0.07916666666666666297
$ws.Range("B112").Value = "Middleware para procesar JSON en el cuerpo de las solicitudes POST  app.use(express.json());"
$ws.Range("A9:B9").Copy() | Out-Null
$ws.Range("A112:B112").PasteSpecial(-4122) | Out-Null

# --- Hyperlink on B105 (https://expressjs.com) ---
$ws.Hyperlinks.Add($ws.Range("B105"), "https://expressjs.com") | Out-Null

# --- Restore intended cell styles for B105 (hyperlink) in case Add() changed it ---
$ws.Range("A21:B21").Copy() | Out-Null
$ws.Range("A105:B105").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(105).RowHeight = 15

# --- Final view state: selection + scroll position ---
$ws.Range("A113").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 79
